# Atividade 1 - add "correct answer" column (AC) for the matching exercise,
# fix one missing score value (Q19), fill in the two summary cells (Z28/AA28),
# and append three new rows of notes at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- New "correct answer" column (AC), populated in the same order the
#     original author typed them so the shared-string table comes out in the
#     same sequence. ---
$ws.Range("AC10").Value = "resposta certa: internet"
$ws.Range("AC12").Value = "resposta certa: Web"
$ws.Range("AC14").Value = "resposta certa: WEB 1.0"
$ws.Range("AC5").Value  = "CERTO"
$ws.Range("AC7").Value  = "CERTO"
$ws.Range("AC8").Value  = "CERTO"
$ws.Range("AC9").Value  = "CERTO"
$ws.Range("AC11").Value = "CERTO"
$ws.Range("AC13").Value = "CERTO"
$ws.Range("AC19").Value = "CERTO"
$ws.Range("AC22").Value = "CERTO"
$ws.Range("AC15").Value = "resposta certa: WEB 2.0"
$ws.Range("AC16").Value = "resposta certa: WEB 3.0"
$ws.Range("AC17").Value = "resposta certa: WEB 4.0"
$ws.Range("AC18").Value = "resposta certa: Internet das coisas"
$ws.Range("AC23").Value = "resposta certa: Internet das coisas"
$ws.Range("AC20").Value = "resposta certa: WEB"
$ws.Range("AC21").Value = "resposta certa: WEB"
$ws.Range("AC24").Value = "resposta certa: Cloud Computing"
$ws.Range("AC26").Value = "resposta certa: Edge Computing"
$ws.Range("AC25").Value = "resposta certa: Fog Computing"

# --- Missing score for "Netscape / 1993" item (row 19), column Q (trio score) ---
$ws.Range("Q19").Value = 12

# --- Fill the two leftover summary cells below the grand totals row ---
$ws.Range("Z28").Value = 139
$ws.Range("AA28").Value = 140

# --- Size the new AC column to fit its content (author's "bestFit" column) ---
$ws.Columns.Item(29).ColumnWidth = 33.8

# --- Three extra note rows at the bottom of the sheet ---
$ws.Range("A32").Value = "Armazenamento"
$ws.Range("A33").Value = "Padronização"
$ws.Range("A34").Value = "Eficiencia Energetica"

# --- Move the selection / scrolled view to where the author left off ---
$ws.Range("D29").Select()
